$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.626.48'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.846.65'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.07'
$ws.Range('E5').Value = '  +1.02%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4316'
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3696'
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.29'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07339'
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8764'
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.05'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '1.879.29'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.477'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.606'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '81.57'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009079'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.59'
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('D22').Value = '27.699.58'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.093'
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.00'
$ws.Range('E24').Value = '  +5.75%  '
$ws.Range('D25').Value = '2.102.87'
$ws.Range('E25').Value = '  -0.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.988'
$ws.Range('E26').Value = '  +1.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.38'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.05'
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.337'
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.99'
$ws.Range('E30').Value = '  -4.61%  '
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08915'
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7860'
$ws.Range('E33').Value = '  +3.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.621'
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.172'
$ws.Range('E35').Value = '  +6.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.983'
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.111'
$ws.Range('E38').Value = '  +1.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05445'
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01967'
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5189'
$ws.Range('E42').Value = '  +2.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1695'
$ws.Range('E43').Value = '  +2.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.776'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.661'
$ws.Range('E45').Value = '  +3.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.62'
$ws.Range('E46').Value = '  +3.14%  '
$ws.Range('E47').Value = '  +2.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.77'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06544'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.005'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.670'
$ws.Range('E51').Value = '  +3.09%  '
